# 012 Week 7 data update
#
# THURSDAY SINGLES: enter the Week 7 ("WK 7", column H) scores that were
# collected for this round. The Z column totals are SUM() formulas and
# recalculate automatically.
#
# HANDICAPS: refresh the handicap figures for a few players following the
# week's results.

$wb = $excel.ActiveWorkbook

# --- THURSDAY SINGLES: Week 7 (column H) scores ---
$thu = $wb.Worksheets.Item("THURSDAY SINGLES")
$thu.Range("H5").Value  = 27
$thu.Range("H6").Value  = 27
$thu.Range("H7").Value  = 26
$thu.Range("H8").Value  = 36
$thu.Range("H9").Value  = 20
$thu.Range("H10").Value = 16
$thu.Range("H11").Value = 23
$thu.Range("H12").Value = 28
$thu.Range("H13").Value = 31
$thu.Range("H15").Value = 34

# --- HANDICAPS: updated handicap values ---
$hc = $wb.Worksheets.Item("HANDICAPS")
$hc.Range("B2").Value  = 12
$hc.Range("C2").Value  = 12
$hc.Range("B6").Value  = 17
$hc.Range("C6").Value  = 17
$hc.Range("C22").Value = 13
$hc.Range("C23").Value = 12
$hc.Range("C24").Value = 14
